# Data-cleanup pass over the "purpose" (C), "age" (F) and "job" (G) columns.
#
# Rules observed in the target diff:
#   Column C (purpose):
#     'new car'              -> new car              (strip wrapping quotes)
#     'used car'             -> used car              (strip wrapping quotes)
#     'domestic appliance'   -> domestic appliance     (strip wrapping quotes)
#     new car'               -> new car               (typo: stray trailing quote)
#     use car'                -> used car              (typo: stray trailing quote + missing "d")
#     the                     -> other                 (bad value)
#   Column F (age):
#     Old                    -> Senior
#   Column G (job):
#     unskilled resident     -> 'unskilled resident'   (missing wrapping quotes)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {

    # --- Column C: purpose -------------------------------------------------
    $c = $ws.Cells.Item($r, 3)
    $cv = $c.Value2

    if ($cv -eq "'new car'") {
        $c.Value = "new car"
    } elseif ($cv -eq "'used car'") {
        $c.Value = "used car"
    } elseif ($cv -eq "'domestic appliance'") {
        $c.Value = "domestic appliance"
    } elseif ($cv -eq "new car'") {
        $c.Value = "new car"
    } elseif ($cv -eq "use car'") {
        $c.Value = "used car"
    } elseif ($cv -eq "the") {
        $c.Value = "other"
    }

    # --- Column F: age -------------------------------------------------------
    $f = $ws.Cells.Item($r, 6)
    $fv = $f.Value2

    if ($fv -eq "Old") {
        $f.Value = "Senior"
    }

    # --- Column G: job -------------------------------------------------------
    $g = $ws.Cells.Item($r, 7)
    $gv = $g.Value2

    if ($gv -eq "unskilled resident") {
        # NB: a single leading apostrophe is treated by Excel as a text-prefix
        # marker and is not stored as literal data, so it has to be doubled
        # here to end up with one literal leading quote in the cell text.
        $g.Value = "''unskilled resident'"
    }
}
